# Auto-generated edit script applying market-price / profit updates
# described by the commit diff, grouped per worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64 (ALC)
$ws.Cells.Item(64, 8).Value = 5811.1113
$ws.Cells.Item(64, 9).Value = 7033.3335
$ws.Cells.Item(64, 10).Value = 3366.6667
$ws.Cells.Item(64, 11).Value = 7033.3335
$ws.Cells.Item(64, 12).Value = 3366.6667
$ws.Cells.Item(64, 13).Value = -6785.3335
$ws.Cells.Item(64, 14).Value = -3862.6667

# Row 67 (ALC)
$ws.Cells.Item(67, 8).Value = 5811.1113
$ws.Cells.Item(67, 9).Value = 7033.3335
$ws.Cells.Item(67, 10).Value = 3366.6667
$ws.Cells.Item(67, 11).Value = 7033.3335
$ws.Cells.Item(67, 12).Value = 3366.6667
$ws.Cells.Item(67, 13).Value = -6175.3335
$ws.Cells.Item(67, 14).Value = -5082.6667

# Row 74 (ALC)
$ws.Cells.Item(74, 8).Value = 4329526
$ws.Cells.Item(74, 9).Value = 6492390
$ws.Cells.Item(74, 10).Value = 3797.5
$ws.Cells.Item(74, 11).Value = 6492390
$ws.Cells.Item(74, 12).Value = 3797.5
$ws.Cells.Item(74, 13).Value = -6491454
$ws.Cells.Item(74, 14).Value = -5669.5

# Row 76 (ALC)
$ws.Cells.Item(76, 8).Value = 53574276
$ws.Cells.Item(76, 9).Value = 60002830
$ws.Cells.Item(76, 10).Value = 2966.6667
$ws.Cells.Item(76, 11).Value = 60002830
$ws.Cells.Item(76, 12).Value = 2966.6667
$ws.Cells.Item(76, 13).Value = -60002515
$ws.Cells.Item(76, 14).Value = -3596.6667

# Row 77 (ALC)
$ws.Cells.Item(77, 8).Value = 4329526
$ws.Cells.Item(77, 9).Value = 6492390
$ws.Cells.Item(77, 10).Value = 3797.5
$ws.Cells.Item(77, 11).Value = 32461950
$ws.Cells.Item(77, 12).Value = 18987.5
$ws.Cells.Item(77, 13).Value = -32457270
$ws.Cells.Item(77, 14).Value = -28347.5

# Row 79 (ALC)
$ws.Cells.Item(79, 8).Value = 53574276
$ws.Cells.Item(79, 9).Value = 60002830
$ws.Cells.Item(79, 10).Value = 2966.6667
$ws.Cells.Item(79, 11).Value = 60002830
$ws.Cells.Item(79, 12).Value = 2966.6667
$ws.Cells.Item(79, 13).Value = -60001738
$ws.Cells.Item(79, 14).Value = -5150.6667

# Row 132 (ALC)
$ws.Cells.Item(132, 8).Value = 1819846.6
$ws.Cells.Item(132, 9).Value = 1671.5116
$ws.Cells.Item(132, 10).Value = 8334974
$ws.Cells.Item(132, 11).Value = 5014.5348
$ws.Cells.Item(132, 12).Value = 25004922
$ws.Cells.Item(132, 13).Value = -2484.5348
$ws.Cells.Item(132, 14).Value = -25009982

# Row 138 (ALC)
$ws.Cells.Item(138, 8).Value = 3770.4092
$ws.Cells.Item(138, 9).Value = 1794.9678
$ws.Cells.Item(138, 10).Value = 4844.772
$ws.Cells.Item(138, 11).Value = 5384.903399999999
$ws.Cells.Item(138, 12).Value = 14534.316
$ws.Cells.Item(138, 13).Value = -244.9033999999992
$ws.Cells.Item(138, 14).Value = -24814.316

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Cells.Item(32, 8).Value = 20098.166
$ws.Cells.Item(32, 9).Value = 15010.2
$ws.Cells.Item(32, 10).Value = 60801.9
$ws.Cells.Item(32, 11).Value = 15010.2
$ws.Cells.Item(32, 12).Value = 60801.9
$ws.Cells.Item(32, 13).Value = -14723.2
$ws.Cells.Item(32, 14).Value = -61375.9

# Row 110 (ARM)
$ws.Cells.Item(110, 8).Value = 714.8889
$ws.Cells.Item(110, 9).Value = 651.1177
$ws.Cells.Item(110, 10).Value = 1799
$ws.Cells.Item(110, 11).Value = 651.1177
$ws.Cells.Item(110, 12).Value = 1799
$ws.Cells.Item(110, 13).Value = 1393.8823
$ws.Cells.Item(110, 14).Value = -5889

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (BSM)
$ws.Cells.Item(105, 8).Value = 3190.625
$ws.Cells.Item(105, 9).Value = 3299.577
$ws.Cells.Item(105, 10).Value = 2718.5
$ws.Cells.Item(105, 11).Value = 3299.577
$ws.Cells.Item(105, 12).Value = 2718.5
$ws.Cells.Item(105, 13).Value = -1552.577
$ws.Cells.Item(105, 14).Value = -6212.5

# Row 141 (BSM)
$ws.Cells.Item(141, 8).Value = 54835.1
$ws.Cells.Item(141, 10).Value = 54835.1
$ws.Cells.Item(141, 12).Value = 54835.1
$ws.Cells.Item(141, 14).Value = -65195.1

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Cells.Item(31, 8).Value = 2889.7273
$ws.Cells.Item(31, 9).Value = 2186.4
$ws.Cells.Item(31, 10).Value = 4120.55
$ws.Cells.Item(31, 11).Value = 2186.4
$ws.Cells.Item(31, 12).Value = 4120.55
$ws.Cells.Item(31, 13).Value = -1891.4
$ws.Cells.Item(31, 14).Value = -4710.55

# Row 34 (CRP)
$ws.Cells.Item(34, 8).Value = 2889.7273
$ws.Cells.Item(34, 9).Value = 2186.4
$ws.Cells.Item(34, 10).Value = 4120.55
$ws.Cells.Item(34, 11).Value = 2186.4
$ws.Cells.Item(34, 12).Value = 4120.55
$ws.Cells.Item(34, 13).Value = -1984.4
$ws.Cells.Item(34, 14).Value = -4524.55

# Row 62 (CRP)
$ws.Cells.Item(62, 8).Value = 5763.9287
$ws.Cells.Item(62, 9).Value = 5789.5
$ws.Cells.Item(62, 10).Value = 5700
$ws.Cells.Item(62, 11).Value = 5789.5
$ws.Cells.Item(62, 12).Value = 5700
$ws.Cells.Item(62, 13).Value = -5165.5
$ws.Cells.Item(62, 14).Value = -6948

# Row 65 (CRP)
$ws.Cells.Item(65, 8).Value = 5763.9287
$ws.Cells.Item(65, 9).Value = 5789.5
$ws.Cells.Item(65, 10).Value = 5700
$ws.Cells.Item(65, 11).Value = 28947.5
$ws.Cells.Item(65, 12).Value = 28500
$ws.Cells.Item(65, 13).Value = -25827.5
$ws.Cells.Item(65, 14).Value = -34740

$ws = $wb.Worksheets.Item("CUL")
# Row 122 (CUL)
$ws.Cells.Item(122, 8).Value = 1502.7941
$ws.Cells.Item(122, 9).Value = 574.25
$ws.Cells.Item(122, 10).Value = 1626.6
$ws.Cells.Item(122, 11).Value = 5168.25
$ws.Cells.Item(122, 12).Value = 14639.4
$ws.Cells.Item(122, 13).Value = -2718.25
$ws.Cells.Item(122, 14).Value = -19539.4

$ws = $wb.Worksheets.Item("GSM")
# Row 70 (GSM)
$ws.Cells.Item(70, 8).Value = 14501
$ws.Cells.Item(70, 9).Value = 19839.2
$ws.Cells.Item(70, 10).Value = 4235.231
$ws.Cells.Item(70, 11).Value = 19839.2
$ws.Cells.Item(70, 12).Value = 4235.231
$ws.Cells.Item(70, 13).Value = -19569.2
$ws.Cells.Item(70, 14).Value = -4775.231

# Row 73 (GSM)
$ws.Cells.Item(73, 8).Value = 14501
$ws.Cells.Item(73, 9).Value = 19839.2
$ws.Cells.Item(73, 10).Value = 4235.231
$ws.Cells.Item(73, 11).Value = 19839.2
$ws.Cells.Item(73, 12).Value = 4235.231
$ws.Cells.Item(73, 13).Value = -18903.2
$ws.Cells.Item(73, 14).Value = -6107.231

# Row 80 (GSM)
$ws.Cells.Item(80, 8).Value = 4030.697
$ws.Cells.Item(80, 9).Value = 4118.0527
$ws.Cells.Item(80, 10).Value = 3912.1428
$ws.Cells.Item(80, 11).Value = 4118.0527
$ws.Cells.Item(80, 12).Value = 3912.1428
$ws.Cells.Item(80, 13).Value = -3120.0527
$ws.Cells.Item(80, 14).Value = -5908.1428

# Row 83 (GSM)
$ws.Cells.Item(83, 8).Value = 4030.697
$ws.Cells.Item(83, 9).Value = 4118.0527
$ws.Cells.Item(83, 10).Value = 3912.1428
$ws.Cells.Item(83, 11).Value = 20590.2635
$ws.Cells.Item(83, 12).Value = 19560.714
$ws.Cells.Item(83, 13).Value = -15598.2635
$ws.Cells.Item(83, 14).Value = -29544.714

# Row 122 (GSM)
$ws.Cells.Item(122, 8).Value = 775662.7
$ws.Cells.Item(122, 9).Value = 1013382.1
$ws.Cells.Item(122, 10).Value = 3074.5
$ws.Cells.Item(122, 11).Value = 3040146.3
$ws.Cells.Item(122, 12).Value = 9223.5
$ws.Cells.Item(122, 13).Value = -3037696.3
$ws.Cells.Item(122, 14).Value = -14123.5

$ws = $wb.Worksheets.Item("LTW")
# Row 30 (LTW)
$ws.Cells.Item(30, 8).Value = 400
$ws.Cells.Item(30, 9).Value = 400
$ws.Cells.Item(30, 10).Value = 0
$ws.Cells.Item(30, 11).Value = 400
$ws.Cells.Item(30, 12).Value = 0
$ws.Cells.Item(30, 13).Value = -292
$ws.Cells.Item(30, 14).ClearContents()

# Row 132 (LTW)
$ws.Cells.Item(132, 8).Value = 2744.0894
$ws.Cells.Item(132, 9).Value = 1656.8462
$ws.Cells.Item(132, 10).Value = 5238.353
$ws.Cells.Item(132, 11).Value = 4970.5386
$ws.Cells.Item(132, 12).Value = 15715.059
$ws.Cells.Item(132, 13).Value = -2440.5386
$ws.Cells.Item(132, 14).Value = -20775.059

# Row 136 (LTW)
$ws.Cells.Item(136, 8).Value = 4138.364
$ws.Cells.Item(136, 9).Value = 2054.5833
$ws.Cells.Item(136, 10).Value = 5751.613
$ws.Cells.Item(136, 11).Value = 6163.749899999999
$ws.Cells.Item(136, 12).Value = 17254.839
$ws.Cells.Item(136, 13).Value = -3613.749899999999
$ws.Cells.Item(136, 14).Value = -22354.839

$ws = $wb.Worksheets.Item("WVR")
# Row 96 (WVR)
$ws.Cells.Item(96, 8).Value = 3324.9167
$ws.Cells.Item(96, 9).Value = 2111.111
$ws.Cells.Item(96, 10).Value = 4053.2
$ws.Cells.Item(96, 11).Value = 2111.111
$ws.Cells.Item(96, 12).Value = 4053.2
$ws.Cells.Item(96, 13).Value = -738.1109999999999
$ws.Cells.Item(96, 14).Value = -6799.2
